$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "28.321.80", "  -5.32%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.840.30", "  -4.99%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.002", "  -0.49%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "329.89", "  -1.82%  "),
    @(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.001", "  -0.40%  "),
    @(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.4602", "  -4.79%  "),
    @(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.3864", "  -6.16%  "),
    @(9, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "46.08", "  -2.68%  "),
    @(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.07868", "  -3.65%  "),
    @(11, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.9656", "  -4.90%  "),
    @(12, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "21.99", "  -7.08%  "),
    @(13, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.862.42", "  -5.42%  "),
    @(14, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.724", "  -5.98%  "),
    @(15, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "6.927", "  -4.99%  "),
    @(16, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.06888", "  +0.62%  "),
    @(17, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.002", "  -0.46%  "),
    @(18, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "87.00", "  -4.57%  "),
    @(19, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000009949", "  -3.92%  "),
    @(20, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "16.95", "  -4.78%  "),
    @(21, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.002", "  -0.33%  "),
    @(22, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "28.344.33", "  -5.22%  "),
    @(23, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "5.343", "  -5.12%  "),
    @(24, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "11.00", "  -7.36%  "),
    @(25, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.150", "  -1.36%  "),
    @(26, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.080.11", "  -5.31%  "),
    @(27, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "153.69", "  -1.99%  "),
    @(28, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "19.24", "  -3.87%  "),
    @(29, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "5.779", "  -13.71%  "),
    @(30, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.990", "  -5.08%  "),
    @(31, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "117.11", "  -3.70%  "),
    @(32, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.9435", "  -6.32%  "),
    @(33, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.09311", "  -3.34%  "),
    @(34, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "5.289", "  -5.19%  "),
    @(35, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "3.443", "  -2.86%  "),
    @(36, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.329", "  -6.40%  "),
    @(37, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.06022", "  -8.14%  "),
    @(38, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.02156", "  -5.67%  "),
    @(39, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.150", "  -4.68%  "),
    @(40, "Frax", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", "1.001", "  -0.39%  "),
    @(41, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "7.624", "  -4.47%  "),
    @(42, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.5635", "  -5.58%  "),
    @(43, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "10.03", "  -6.49%  "),
    @(44, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1783", "  -3.45%  "),
    @(45, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "2.282", "  -8.38%  "),
    @(46, "WEMIXToken", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "1.213", "  -4.38%  "),
    @(47, "Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "0.5299", "  -4.86%  "),
    @(48, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "11.62", "  -6.18%  "),
    @(49, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.07039", "  -5.92%  "),
    @(50, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "1.843", "  -7.38%  "),
    @(51, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "112.90", "  -3.54%  ")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
